# Update NMA and MA coefficient tables to use M instead of mu, for
# consistency with the PDF documentation.
#
# Every "mu_<n>" label in column H (the mu_name / coefficient-name
# column) becomes "M_<n>" on all four sheets: weibull, gompertz,
# fracpoly1, fracpoly2. Header cells such as "mu_name" are left alone.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rows = $used.Rows.Count
    for ($r = 1; $r -le $rows; $r++) {
        $cell = $ws.Cells.Item($r, 8)  # column H
        $s = $cell.Text
        if ($s -match "^mu_\d+$") {
            $cell.Value = "M_" + $s.Substring(3)
        }
    }
}

# Restore/update each sheet's selection.
$wsWeibull = $wb.Worksheets.Item("weibull")
$wsWeibull.Activate()
$wsWeibull.Range("H8").Select()

$wsGompertz = $wb.Worksheets.Item("gompertz")
$wsGompertz.Activate()
$wsGompertz.Range("H8").Select()

$wsFracpoly1 = $wb.Worksheets.Item("fracpoly1")
$wsFracpoly1.Activate()
$wsFracpoly1.Range("H10").Select()

$wsFracpoly2 = $wb.Worksheets.Item("fracpoly2")
$wsFracpoly2.Activate()
$wsFracpoly2.Range("H5").Select()

# fracpoly2 is the sheet left active/selected in the workbook.
$wsFracpoly2.Activate()
